# Add the new "dev.char" worksheet right after "head_width"
$wb = $excel.ActiveWorkbook
$headWidth = $wb.Worksheets.Item("head_width")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $headWidth)
$ws.Name = "dev.char"

# ---- Cell values -----------------------------------------------------
# NOTE: the write order below intentionally mirrors how the shared-string
# table was populated in the source workbook (first three header cells,
# then the table body, then the two extra header cells) so the resulting
# <sst> unique-string order lines up with the authored file.
$ws.Range("A1").Value = "Stage"
$ws.Range("B1").Value = "SET"
$ws.Range("C1").Value = "std.error"

$ws.Range("A2").Value = "Egg"
$ws.Range("B2").Value = 929.35379999999998
$ws.Range("C2").Value = 49.110779999999998

$ws.Range("A3").Value = "L1"
$ws.Range("B3").Value = 233.68279999999999
$ws.Range("C3").Value = 27.030670000000001

$ws.Range("A4").Value = "L2"
$ws.Range("B4").Value = 243.9451
$ws.Range("C4").Value = 45.301189999999998

$ws.Range("A5").Value = "L3"
$ws.Range("B5").Value = 2602.9964
$ws.Range("C5").Value = 297.46382

$ws.Range("A6").Value = "Pupae"
$ws.Range("B6").Value = 1207.4305999999999
$ws.Range("C6").Value = 489.28820999999999

$ws.Range("D1").Value = "Tmin"
$ws.Range("E1").Value = "std.error"

$ws.Range("D2").Value = 11.399692
$ws.Range("E2").Value = 0.36807020000000001

$ws.Range("D3").Value = 15.436838
$ws.Range("E3").Value = 0.30478480000000002

$ws.Range("D4").Value = 15.689105
$ws.Range("E4").Value = 0.40978490000000001

$ws.Range("D5").Value = 9.3754659999999994
$ws.Range("E5").Value = 0.84568220000000005

$ws.Range("D6").Value = 12.535199
$ws.Range("E6").Value = 1.6237360999999999

# Row 7 stays blank (trailing empty row like the original sheet had)

# ---- Number formatting -------------------------------------------------
$ws.Range("B2:E6").NumberFormat = "0.000"

# ---- Fonts ---------------------------------------------------------
# Bold labels in column A (Stage names)
$ws.Range("A1:A6").Font.Bold = $true

# Header row (B1:E1) also bold
$ws.Range("B1:E1").Font.Bold = $true

# Numeric "SET" column (B) - plain Calibri
$ws.Range("B2:B6").Font.Name = "Calibri"
$ws.Range("B2:B6").Font.Size = 11

# std.error / Tmin columns (C,D,E) - Lucida Console, size 10, vertical centred
$ws.Range("C2:E6").Font.Name = "Lucida Console"
$ws.Range("C2:E6").Font.Size = 10
$ws.Range("C2:E6").VerticalAlignment = -4108

# ---- Borders ---------------------------------------------------------
# Column B gets a left border down the whole table
$ws.Range("B1:B6").Borders.Item(7).LineStyle = 1
$ws.Range("B1:B6").Borders.Item(7).Weight = 2

# A horizontal line under the header row
$ws.Range("A2:E2").Borders.Item(8).LineStyle = 1
$ws.Range("A2:E2").Borders.Item(8).Weight = 2

# A horizontal line closing the table under the last data row
$ws.Range("A6:E6").Borders.Item(9).LineStyle = 1
$ws.Range("A6:E6").Borders.Item(9).Weight = 2

# Column sizing to match source
$ws.Columns.Item(1).ColumnWidth = 9.140625

# ---- View / selection -------------------------------------------------
$ws.Range("G5").Select()
